# "New benchmark update - implemented new stage distribution"
# Insert three new rows at the top of the age-group table for the new
# childhood age bins (0,4], (4,9], (9,14] with zero rates, pushing the
# existing 15-19 ... 85+ rows (and the chart that reads them) down by
# three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (rows 2:16) down to rows 5:19.
$ws.Rows("2:4").Insert()

# New rows: (0,4], (4,9], (9,14] — all rates are 0.
$ws.Range("A2").Value = "(0,4]"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("A3").Value = "(4,9]"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("A4").Value = "(9,14]"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

# Match the selection left behind in the saved workbook.
$null = $ws.Range("A2:D4").Select()

# Re-point the chart series at the shifted data range (A5:A19 etc.) and
# move the chart itself down by three rows (3 x 15pt default row height)
# so it still sits below the table.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.SeriesCollection(1).Formula = "=SERIES(""Male and female"",'Sheet 1'!`$A`$5:`$A`$19,'Sheet 1'!`$B`$5:`$B`$19,1)"
$chart.SeriesCollection(2).Formula = "=SERIES(""Male"",'Sheet 1'!`$A`$5:`$A`$19,'Sheet 1'!`$C`$5:`$C`$19,2)"
$chart.SeriesCollection(3).Formula = "=SERIES(""Female"",'Sheet 1'!`$A`$5:`$A`$19,'Sheet 1'!`$D`$5:`$D`$19,3)"

$co.Top = $co.Top + 45
